# Regenerate save_data to use K (Strikes) instead of Strike# for the K column.
# This rewrites column G (header "K") for rows 2-23 with the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 3
    8  = 1
    9  = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 2
    14 = 0
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 2
    23 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
